# Generate Report for Handoff
# Update the "Latest Handoff" timestamps for the 515add6f-fa69-4227-9c96-aa97448c4e51.md
# file, reflecting a fresh handoff report run for both target languages (zh-cn, de-de)
# and the rollup Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-29-18 08:29:17"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-18 08:29:14"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-18 08:29:17"
